$d = $word.ActiveDocument

# Update the date heading paragraph
$d.Content.Find.Execute("2025-05-29 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-05-30 Friday", 2)

# Update each math-fact cell in the results table (row-major, 5 columns)
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "73-51=22"
$t.Cell(1,2).Range.Text = "18+74=92"
$t.Cell(1,3).Range.Text = "65-16=49"
$t.Cell(1,4).Range.Text = "35+13=48"
$t.Cell(1,5).Range.Text = "22-11=11"
$t.Cell(2,1).Range.Text = "88-81=7"
$t.Cell(2,2).Range.Text = "53+41=94"
$t.Cell(2,3).Range.Text = "73+21=94"
$t.Cell(2,4).Range.Text = "28+9=37"
$t.Cell(2,5).Range.Text = "97-88=9"
$t.Cell(3,1).Range.Text = "4-4=0"
$t.Cell(3,2).Range.Text = "67-52=15"
$t.Cell(3,3).Range.Text = "67-63=4"
$t.Cell(3,4).Range.Text = "47+36=83"
$t.Cell(3,5).Range.Text = "43+53=96"
$t.Cell(4,1).Range.Text = "77-38=39"
$t.Cell(4,2).Range.Text = "93-82=11"
$t.Cell(4,3).Range.Text = "34+63=97"
$t.Cell(4,4).Range.Text = "81-43=38"
$t.Cell(4,5).Range.Text = "52+26=78"
$t.Cell(5,1).Range.Text = "40+0=40"
$t.Cell(5,2).Range.Text = "95-42=53"
$t.Cell(5,3).Range.Text = "66-26=40"
$t.Cell(5,4).Range.Text = "71-7=64"
$t.Cell(5,5).Range.Text = "42+57=99"
$t.Cell(6,1).Range.Text = "24+73=97"
$t.Cell(6,2).Range.Text = "39+34=73"
$t.Cell(6,3).Range.Text = "70-12=58"
$t.Cell(6,4).Range.Text = "16+50=66"
$t.Cell(6,5).Range.Text = "81+8=89"
$t.Cell(7,1).Range.Text = "51+39=90"
$t.Cell(7,2).Range.Text = "74-4=70"
$t.Cell(7,3).Range.Text = "25+31=56"
$t.Cell(7,4).Range.Text = "1+19=20"
$t.Cell(7,5).Range.Text = "61-6=55"
$t.Cell(8,1).Range.Text = "63+16=79"
$t.Cell(8,2).Range.Text = "79-72=7"
$t.Cell(8,3).Range.Text = "56+42=98"
$t.Cell(8,4).Range.Text = "52-24=28"
$t.Cell(8,5).Range.Text = "97-60=37"
$t.Cell(9,1).Range.Text = "10+40=50"
$t.Cell(9,2).Range.Text = "55-7=48"
$t.Cell(9,3).Range.Text = "91-49=42"
$t.Cell(9,4).Range.Text = "19+44=63"
$t.Cell(9,5).Range.Text = "75-24=51"
$t.Cell(10,1).Range.Text = "15+10=25"
$t.Cell(10,2).Range.Text = "31+20=51"
$t.Cell(10,3).Range.Text = "52+13=65"
$t.Cell(10,4).Range.Text = "59-45=14"
$t.Cell(10,5).Range.Text = "76-59=17"
$t.Cell(11,1).Range.Text = "15-7=8"
$t.Cell(11,2).Range.Text = "14+32=46"
$t.Cell(11,3).Range.Text = "48-47=1"
$t.Cell(11,4).Range.Text = "45+15=60"
$t.Cell(11,5).Range.Text = "48-18=30"
$t.Cell(12,1).Range.Text = "32-5=27"
$t.Cell(12,2).Range.Text = "40+23=63"
$t.Cell(12,3).Range.Text = "1+22=23"
$t.Cell(12,4).Range.Text = "39-38=1"
$t.Cell(12,5).Range.Text = "96-3=93"
$t.Cell(13,1).Range.Text = "25-16=9"
$t.Cell(13,2).Range.Text = "22+34=56"
$t.Cell(13,3).Range.Text = "70+21=91"
$t.Cell(13,4).Range.Text = "31-25=6"
$t.Cell(13,5).Range.Text = "7+73=80"
$t.Cell(14,1).Range.Text = "19+18=37"
$t.Cell(14,2).Range.Text = "84+1=85"
$t.Cell(14,3).Range.Text = "42-7=35"
$t.Cell(14,4).Range.Text = "78-23=55"
$t.Cell(14,5).Range.Text = "91-0=91"
$t.Cell(15,1).Range.Text = "58+17=75"
$t.Cell(15,2).Range.Text = "29+3=32"
$t.Cell(15,3).Range.Text = "0+55=55"
$t.Cell(15,4).Range.Text = "43+56=99"
$t.Cell(15,5).Range.Text = "41+1=42"
$t.Cell(16,1).Range.Text = "10+20=30"
$t.Cell(16,2).Range.Text = "78-73=5"
$t.Cell(16,3).Range.Text = "32+26=58"
$t.Cell(16,4).Range.Text = "54-32=22"
$t.Cell(16,5).Range.Text = "71-27=44"
$t.Cell(17,1).Range.Text = "3+42=45"
$t.Cell(17,2).Range.Text = "8-0=8"
$t.Cell(17,3).Range.Text = "8+1=9"
$t.Cell(17,4).Range.Text = "80-54=26"
$t.Cell(17,5).Range.Text = "64-34=30"
$t.Cell(18,1).Range.Text = "14+55=69"
$t.Cell(18,2).Range.Text = "11+65=76"
$t.Cell(18,3).Range.Text = "0+30=30"
$t.Cell(18,4).Range.Text = "46+20=66"
$t.Cell(18,5).Range.Text = "41+15=56"
$t.Cell(19,1).Range.Text = "88-67=21"
$t.Cell(19,2).Range.Text = "49-35=14"
$t.Cell(19,3).Range.Text = "55-23=32"
$t.Cell(19,4).Range.Text = "28-10=18"
$t.Cell(19,5).Range.Text = "2+21=23"
$t.Cell(20,1).Range.Text = "63-24=39"
$t.Cell(20,2).Range.Text = "78-35=43"
$t.Cell(20,3).Range.Text = "21-15=6"
$t.Cell(20,4).Range.Text = "27+45=72"
$t.Cell(20,5).Range.Text = "82-53=29"

Write-Output "done"
